# Version 1.1.8 - add "ESTADO" column when exporting.
# The template gains a new column (AS) right before the old
# "COMENTARIOS FACTURA" column, shifting everything from the old AS
# column onward one place to the right, and the new column gets the
# header "ESTADO".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column 45 is "AS" - insert a new blank column there, pushing the
# existing AS:AX columns (and their data) one column to the right.
$ws.Columns(45).Insert()

# New AS8 header cell.
$ws.Cells.Item(8, 45).Value = "ESTADO"

# Reflect the new selection position (matches the author ending up with
# the newly added header cell selected).
$ws.Range("AS8").Select()
